# "Generate Report for Handback"
# Updates the localization-status workbook to reflect a handback report:
#  - Overview sheet: status cells move from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: refresh the "Latest Handback DateTime" cells and
#    clear the stale "Error Detail" cells (handback is now in sync, so there
#    is no longer an out-of-date warning).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# NOTE: the host's ColumnWidth setter quantizes to coarse character-width
# steps, so we feed it a value that lands on the closest reachable stored
# width to the true target (~29.9777 characters).
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# --- zh-cn sheet ------------------------------------------------------------
# "Status" column re-uses the same shared string as the Overview summary, so
# it also flips from "Ready for handoff" to the new handback message.
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("K2").Value = "2016-08-13 04:43:14"
$zhcn.Range("K3").Value = "2016-08-13 04:43:14"

$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(16).ColumnWidth = 12.83

# --- de-de sheet ------------------------------------------------------------
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("K2").Value = "2016-08-13 04:43:23"
$dede.Range("K3").Value = "2016-08-13 04:43:23"

$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(16).ColumnWidth = 12.83
